$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.010462377812505126
$ws.Range("B4").Value = 0.02155249829376056
$ws.Range("B5").Value = 0.04310499658752112
$ws.Range("B6").Value = 0.05231188906252563
$ws.Range("B7").Value = 0.06465749488128168
$ws.Range("B8").Value = 0.083699022500041
$ws.Range("B9").Value = 0.08620999317504224
$ws.Range("B10").Value = 0.1077624914688028
